# Gate_Closure_Trigger.xlsx update
# ---------------------------------------------------------------------
# Add MS River gauge stations (WestPoint, Alliance, Carrollton, BCSpillway,
# BCSpillwayN, Reserve) as new rows 23-28, and the Lower Atchafalaya River
# "MorganCity" station as row 29 -- switching the NOV-14 observation source
# from USACE to USGS.
#
# Column A: short station code. New USGS-style numeric format (0.00000000)
#           with a 10pt Arial font, distinct from the legacy text-format
#           station codes above it.
# Column B: long/descriptive station name (default formatting).
# Column C: observed value, using the existing "0.0" numeric style.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Short codes / long names for the six MS River stations (rows 23-28).
$msRiverCodes = @("WestPoint", "Alliance", "Carrollton", "BCSpillway", "BCSpillwayN", "Reserve")
$msRiverNames = @(
    "MS River at West Point a la Hache",
    "MS River at Alliance",
    "MS River at Carrollton",
    "MS River at Bonnet Carre Spillway",
    "MS River at Bonnet Carre Spillway North",
    "MS River at Reserve"
)

# ---------------------------------------------------------------------
# Column A, rows 23-28: build the new cell format once on A23 (number
# format + Arial 10 black font), then propagate that format to A24:A28
# via Copy (so we don't mint a brand-new font for every row), overwriting
# the value right after each copy with the correct station code.
# ---------------------------------------------------------------------
$a23 = $ws.Range("A23")
$a23.Value = $msRiverCodes[0]
$a23.NumberFormat = "0.00000000"
$fnt = $a23.Font
$fnt.Name = "Arial"
$fnt.Size = 10
$fnt.Color = 0

for ($i = 1; $i -lt $msRiverCodes.Length; $i++) {
    $row = 23 + $i
    $cell = $ws.Range("A" + $row)
    $a23.Copy($cell)
    $cell.Value = $msRiverCodes[$i]
}

# ---------------------------------------------------------------------
# Column B, rows 23-28: plain descriptive names, default formatting.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $msRiverNames.Length; $i++) {
    $row = 23 + $i
    $ws.Range("B" + $row).Value = $msRiverNames[$i]
}

# ---------------------------------------------------------------------
# Column C, rows 23-28: observed values (all 0 for the new stations),
# reusing the workbook's existing "0.0" numeric style.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $row = 23 + $i
    $cCell = $ws.Range("C" + $row)
    $cCell.NumberFormat = "0.0"
    $cCell.Value = 0
}

# ---------------------------------------------------------------------
# Row 29: MorganCity / Lower Atchafalaya River at Morgan City.
# ---------------------------------------------------------------------
$a29 = $ws.Range("A29")
$a23.Copy($a29)
$a29.Value = "MorganCity"

$ws.Range("B29").Value = "Lower Atchafalaya River at Morgan City"

$c29 = $ws.Range("C29")
$c29.NumberFormat = "0.0"
$c29.Value = 0

$ws.Range("A29:C29").Select() | Out-Null

Write-Output "Added MS River stations (rows 23-28) and MorganCity (row 29)"
